# Portfolio Optimisation workbook update
# ----------------------------------------
# "some fixes and ticker selection part for portfolio optimization"
#
#  - Append the newly-closed monthly return rows (Aug-2024 / Sep-2024,
#    Excel serials 45505 / 45536) to the per-ticker return sheets:
#      Sheet2 -> Date, IPB, IIH, TLE
#      Sheet3 -> Date, IIH, TLE
#    (Sheet1, the combined BTC/IPB/IIH/TLE sheet, already has all of its
#    rows; only its selection moved.)
#  - Leave the cursor/selection where the user was last working on each
#    sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws3 = $wb.Worksheets.Item("Sheet3")

# --- Sheet2: columns A (Date), B (IPB), C (IIH), D (TLE) --------------
$ws2.Range("A57").Value = 45505
$ws2.Range("B57").Value = -1.637
$ws2.Range("C57").Value = -5.419
$ws2.Range("D57").Value = 4.744

$ws2.Range("A58").Value = 45536
$ws2.Range("B58").Value = 0.238
$ws2.Range("C58").Value = 2.403
$ws2.Range("D58").Value = 1.641

# Carry the date-column formatting (style used by A2:A56, mmm-yy) onto
# the two new date cells.
$ws2.Range("A56").Copy()
$ws2.Range("A57:A58").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Sheet3: columns A (Date), B (IIH), C (TLE) ------------------------
$ws3.Range("A57").Value = 45505
$ws3.Range("B57").Value = -5.419
$ws3.Range("C57").Value = 4.744

$ws3.Range("A58").Value = 45536
$ws3.Range("B58").Value = 2.403
$ws3.Range("C58").Value = 1.641

$ws3.Range("A56").Copy()
$ws3.Range("A57:A58").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Restore each sheet's working selection / scroll position ---------
$ws1.Range("D6").Select()

$ws2.Activate()
$ws2.Range("A57:D58").Select()
$excel.ActiveWindow.ScrollRow = 49

$ws3.Activate()
$ws3.Range("C61").Select()
$excel.ActiveWindow.ScrollRow = 52
